$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").Value = "`tCategoría"
$ws.Range("B1").Select()
